# ddepewRubric.xlsx update
# - Barrel now uses its proper texture (E63/F63 row "Ability to appropriately translate the view")
#   -> rubric cells in column F (the "X" marker column) get marked for several feature rows,
#      which drives the dependent SUMIF/point formulas to recalculate.
# - Rubric updated for points: mark column F "X" for rows 8, 21, 37, 57, 63, 64, 65.
# - Move the visible selection/viewport to where the grader was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$rowsToMark = @(8, 21, 37, 57, 63, 64, 65)
foreach ($r in $rowsToMark) {
    $ws.Range("F$r").Value = "X"
}

# Update the view: scroll position and active selection, matching the author's
# last on-screen position when they saved the workbook.
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F56").Select() | Out-Null
